$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($range, [string]$text)
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.Style = "Normal"
}

# Per-row Price (D) / Volume(1h) (E) updates
Set-TextValue $ws.Range("D2") "45.992.66"
Set-TextValue $ws.Range("E2") "  -1.79%  "
Set-TextValue $ws.Range("D3") "2.382.27"
Set-TextValue $ws.Range("E3") "  +2.86%  "
Set-TextValue $ws.Range("D4") "1.00"
Set-TextValue $ws.Range("E4") "  +0.20%  "
Set-TextValue $ws.Range("D5") "301.38"
Set-TextValue $ws.Range("E5") "  -0.71%  "
Set-TextValue $ws.Range("D6") "98.87"
Set-TextValue $ws.Range("E6") "  -4.67%  "
Set-TextValue $ws.Range("D7") "0.570"
Set-TextValue $ws.Range("E7") "  -0.82%  "
Set-TextValue $ws.Range("D8") "1.00"
Set-TextValue $ws.Range("E8") "  +0.20%  "
Set-TextValue $ws.Range("D9") "0.514"
Set-TextValue $ws.Range("E9") "  -3.88%  "
Set-TextValue $ws.Range("D10") "34.44"
Set-TextValue $ws.Range("E10") "  -8.27%  "
Set-TextValue $ws.Range("D11") "0.0793"
Set-TextValue $ws.Range("E11") "  -1.97%  "
Set-TextValue $ws.Range("D12") "7.16"
Set-TextValue $ws.Range("E12") "  -4.34%  "
Set-TextValue $ws.Range("E13") "  -0.79%  "
Set-TextValue $ws.Range("D14") "2.749.99"
Set-TextValue $ws.Range("E14") "  +3.16%  "
Set-TextValue $ws.Range("D15") "2.386.66"
Set-TextValue $ws.Range("E15") "  +3.26%  "
Set-TextValue $ws.Range("D16") "0.817"
Set-TextValue $ws.Range("E16") "  -1.52%  "
Set-TextValue $ws.Range("E17") "  -3.10%  "
Set-TextValue $ws.Range("D18") "45.985.08"
Set-TextValue $ws.Range("E18") "  -1.72%  "
Set-TextValue $ws.Range("D19") "12.78"
Set-TextValue $ws.Range("E19") "  -5.34%  "
Set-TextValue $ws.Range("D20") "0.0₃0960"
Set-TextValue $ws.Range("E20") "  +0.47%  "
Set-TextValue $ws.Range("D21") "6.04"
Set-TextValue $ws.Range("E21") "  -2.62%  "
Set-TextValue $ws.Range("D22") "67.57"
Set-TextValue $ws.Range("E22") "  +0.31%  "
Set-TextValue $ws.Range("D23") "244.66"
Set-TextValue $ws.Range("E23") "  -1.93%  "
Set-TextValue $ws.Range("E24") "  -5.15%  "
Set-TextValue $ws.Range("D25") "1.94"
Set-TextValue $ws.Range("E25") "  -3.05%  "
Set-TextValue $ws.Range("E26") "  -0.13%  "
Set-TextValue $ws.Range("D27") "39.61"
Set-TextValue $ws.Range("E27") "  -9.21%  "
Set-TextValue $ws.Range("E28") "  -3.45%  "
Set-TextValue $ws.Range("D29") "9.79"
Set-TextValue $ws.Range("E29") "  -2.46%  "
Set-TextValue $ws.Range("D30") "3.78"
Set-TextValue $ws.Range("E30") "  +19.68%  "
Set-TextValue $ws.Range("D31") "21.13"
Set-TextValue $ws.Range("E31") "  +4.07%  "
Set-TextValue $ws.Range("E32") "  +7.12%  "
Set-TextValue $ws.Range("D33") "5.53"
Set-TextValue $ws.Range("E33") "  -5.20%  "
Set-TextValue $ws.Range("D34") "146.83"
Set-TextValue $ws.Range("E34") "  -0.33%  "
Set-TextValue $ws.Range("D35") "0.0777"
Set-TextValue $ws.Range("E35") "  -3.72%  "
Set-TextValue $ws.Range("E36") "  -0.14%  "
Set-TextValue $ws.Range("D37") "1.93"
Set-TextValue $ws.Range("E37") "  +6.73%  "
Set-TextValue $ws.Range("E38") "  -3.25%  "
Set-TextValue $ws.Range("D39") "14.84"
Set-TextValue $ws.Range("E39") "  -6.83%  "
Set-TextValue $ws.Range("D40") "3.93"
Set-TextValue $ws.Range("E40") "  -3.71%  "
Set-TextValue $ws.Range("D41") "0.0300"
Set-TextValue $ws.Range("E41") "  -3.13%  "
Set-TextValue $ws.Range("D42") "3.21"
Set-TextValue $ws.Range("E42") "  -8.23%  "
Set-TextValue $ws.Range("D43") "1.925.24"
Set-TextValue $ws.Range("E43") "  +3.81%  "
Set-TextValue $ws.Range("E44") "  +0.02%  "
Set-TextValue $ws.Range("D45") "91.53"
Set-TextValue $ws.Range("E45") "  +1.92%  "
Set-TextValue $ws.Range("D46") "1.77"
Set-TextValue $ws.Range("E46") "  -12.62%  "
Set-TextValue $ws.Range("D47") "8.45"
Set-TextValue $ws.Range("E47") "  +5.04%  "
Set-TextValue $ws.Range("D48") "0.186"
Set-TextValue $ws.Range("E48") "  -6.21%  "
Set-TextValue $ws.Range("D51") "68.49"
Set-TextValue $ws.Range("E51") "  -9.22%  "

# Rows 49 and 50: content swapped (Aave <-> RocketPoolETH) with updated prices
Set-TextValue $ws.Range("B49") "RocketPoolETH"
Set-TextValue $ws.Range("C49") "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
Set-TextValue $ws.Range("D49") "2.621.48"
Set-TextValue $ws.Range("E49") "  +3.11%  "

Set-TextValue $ws.Range("B50") "Aave"
Set-TextValue $ws.Range("C50") "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
Set-TextValue $ws.Range("D50") "98.09"
Set-TextValue $ws.Range("E50") "  -0.22%  "
